# Swap columns D (codeforiati:group-name) and E (codeforiati:group-code)
# so that column D becomes codeforiati:group-code and column E becomes
# codeforiati:group-name, for the header row and every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)

    $dVal = $dCell.Value()
    $eVal = $eCell.Value()

    $dCell.Value = $eVal
    $eCell.Value = $dVal
}
